$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Invalid): only row 3 changes from 0 to 1
$ws.Range("G3").Value = 1

# Column H (Absent): rows 3 through 18 change from 0 to 1
$ws.Range("H3:H18").Value = 1
